# Scheduled-runner update: refresh currentAveragePrice / LevePrice* / LeveProfit*
# columns (H, I, J, K, L, M, N) for a batch of leve rows across several
# sheets, per the latest market-board pull. Cells that end up at 0 for every
# tracked column drop their (now meaningless) profit cell entirely; cells
# that newly show a non-zero price pick up a profit cell.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 29999.5
$ws.Range("J93").Value = 29999.5
$ws.Range("L93").Value = 29999.5
$ws.Range("N93").Value = -34991.5

$ws.Range("H95").Value = 29229.6
$ws.Range("J95").Value = 29229.6
$ws.Range("L95").Value = 29229.6
$ws.Range("N95").Value = -34721.6

$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = ""

$ws.Range("H137").Value = 5403.3438
$ws.Range("I137").Value = 5907.9
$ws.Range("K137").Value = 17723.7
$ws.Range("M137").Value = -15173.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3942.817
$ws.Range("I32").Value = 4466.915
$ws.Range("K32").Value = 4466.915
$ws.Range("M32").Value = -4179.915

$ws.Range("H61").Value = 4990.273
$ws.Range("I61").Value = 4611.75
$ws.Range("K61").Value = 4611.75
$ws.Range("M61").Value = -4399.75

$ws.Range("H76").Value = 50288
$ws.Range("J76").Value = 50288
$ws.Range("L76").Value = 50288
$ws.Range("N76").Value = -50964

$ws.Range("H79").Value = 50288
$ws.Range("J79").Value = 50288
$ws.Range("L79").Value = 50288
$ws.Range("N79").Value = -52628

$ws.Range("H95").Value = 65716.836
$ws.Range("J95").Value = 65716.836
$ws.Range("L95").Value = 65716.836
$ws.Range("N95").Value = -71208.836

$ws.Range("H96").Value = 36781.168
$ws.Range("J96").Value = 36781.168
$ws.Range("L96").Value = 36781.168
$ws.Range("N96").Value = -42273.168

$ws.Range("H102").Value = 3402.3333
$ws.Range("I102").Value = 1725.591
$ws.Range("J102").Value = 10780
$ws.Range("K102").Value = 1725.591
$ws.Range("L102").Value = 10780
$ws.Range("M102").Value = -103.5909999999999
$ws.Range("N102").Value = -14024

$ws.Range("H103").Value = 29995
$ws.Range("J103").Value = 29995
$ws.Range("L103").Value = 29995
$ws.Range("N103").Value = -32339

$ws.Range("H122").Value = 3650.6572
$ws.Range("I122").Value = 3233.0344
$ws.Range("J122").Value = 5669.1665
$ws.Range("K122").Value = 9699.1032
$ws.Range("L122").Value = 17007.4995
$ws.Range("M122").Value = -7249.1032
$ws.Range("N122").Value = -21907.4995

$ws.Range("H132").Value = 21278490
$ws.Range("I132").Value = 25642630
$ws.Range("K132").Value = 76927890
$ws.Range("M132").Value = -76925360

$ws.Range("H136").Value = 4990.273
$ws.Range("I136").Value = 4611.75
$ws.Range("K136").Value = 13835.25
$ws.Range("M136").Value = -11285.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5000
$ws.Range("I20").Value = 5000
$ws.Range("K20").Value = 5000
$ws.Range("M20").Value = -4753

$ws.Range("H86").Value = 16668056
$ws.Range("I86").Value = 20834634
$ws.Range("K86").Value = 20834634
$ws.Range("M86").Value = -20833511

$ws.Range("H89").Value = 16668056
$ws.Range("I89").Value = 20834634
$ws.Range("K89").Value = 104173170
$ws.Range("M89").Value = -104167554

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = ""

$ws.Range("H107").Value = 3995.6445
$ws.Range("I107").Value = 3584.7878
$ws.Range("J107").Value = 5125.5
$ws.Range("K107").Value = 3584.7878
$ws.Range("L107").Value = 5125.5
$ws.Range("M107").Value = -1664.7878
$ws.Range("N107").Value = -8965.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = ""

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = ""

$ws.Range("H96").Value = 13587
$ws.Range("J96").Value = 13587
$ws.Range("L96").Value = 13587
$ws.Range("N96").Value = -19079

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 764.9375
$ws.Range("I5").Value = 426.25
$ws.Range("J5").Value = 1781
$ws.Range("K5").Value = 1278.75
$ws.Range("L5").Value = 5343
$ws.Range("M5").Value = -1166.75
$ws.Range("N5").Value = -5567

$ws.Range("H12").Value = 983.37036
$ws.Range("J12").Value = 1038.9546
$ws.Range("L12").Value = 3116.8638
$ws.Range("N12").Value = -3462.8638

$ws.Range("H39").Value = 2821.2632
$ws.Range("J39").Value = 3549.7144
$ws.Range("L39").Value = 10649.1432
$ws.Range("N39").Value = -11237.1432

$ws.Range("H55").Value = 38150296
$ws.Range("I55").Value = 1266.3334
$ws.Range("J55").Value = 57224810
$ws.Range("K55").Value = 3799.0002
$ws.Range("L55").Value = 171674430
$ws.Range("M55").Value = -3622.0002
$ws.Range("N55").Value = -171674784

$ws.Range("H86").Value = 1046.4166
$ws.Range("I86").Value = 1064.4445
$ws.Range("J86").Value = 992.3333
$ws.Range("K86").Value = 3193.3335
$ws.Range("L86").Value = 2976.9999
$ws.Range("M86").Value = -2007.3335
$ws.Range("N86").Value = -5348.9999

$ws.Range("H89").Value = 1046.4166
$ws.Range("I89").Value = 1064.4445
$ws.Range("J89").Value = 992.3333
$ws.Range("K89").Value = 9580.0005
$ws.Range("L89").Value = 8930.9997
$ws.Range("M89").Value = -3652.0005
$ws.Range("N89").Value = -20786.9997

$ws.Range("H132").Value = 3998.75
$ws.Range("I132").Value = 3998
$ws.Range("K132").Value = 35982
$ws.Range("M132").Value = -33452

$ws.Range("H135").Value = 764.9375
$ws.Range("I135").Value = 426.25
$ws.Range("J135").Value = 1781
$ws.Range("K135").Value = 3836.25
$ws.Range("L135").Value = 16029
$ws.Range("M135").Value = -1301.25
$ws.Range("N135").Value = -21099

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 455921.22
$ws.Range("I14").Value = 578961.5600000001
$ws.Range("K14").Value = 578961.5600000001
$ws.Range("M14").Value = -578793.5600000001

$ws.Range("H97").Value = 6263.706
$ws.Range("I97").Value = 436.33334
$ws.Range("K97").Value = 436.33334
$ws.Range("M97").Value = 59.66665999999998

$ws.Range("H101").Value = 31657
$ws.Range("J101").Value = 31657
$ws.Range("L101").Value = 31657
$ws.Range("N101").Value = -38147

$ws.Range("H132").Value = 6633.4053
$ws.Range("I132").Value = 6180.3335
$ws.Range("J132").Value = 7856.7
$ws.Range("K132").Value = 18541.0005
$ws.Range("L132").Value = 23570.1
$ws.Range("M132").Value = -16011.0005
$ws.Range("N132").Value = -28630.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 31250246
$ws.Range("I16").Value = 35714496
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 35714496
$ws.Range("L16").Value = 500
$ws.Range("M16").Value = -35714326
$ws.Range("N16").Value = -840

$ws.Range("H22").Value = 15152968
$ws.Range("J22").Value = 1993.25
$ws.Range("L22").Value = 1993.25
$ws.Range("N22").Value = -2583.25

$ws.Range("H27").Value = 15152968
$ws.Range("J27").Value = 1993.25
$ws.Range("L27").Value = 1993.25
$ws.Range("N27").Value = -2207.25

$ws.Range("H46").Value = 2682.2307
$ws.Range("I46").Value = 1375.25
$ws.Range("K46").Value = 1375.25
$ws.Range("M46").Value = -1187.25

$ws.Range("H55").Value = 335.8
$ws.Range("I55").Value = 362.22223
$ws.Range("J55").Value = 246.625
$ws.Range("K55").Value = 362.22223
$ws.Range("L55").Value = 246.625
$ws.Range("M55").Value = -189.22223
$ws.Range("N55").Value = -592.625

$ws.Range("H61").Value = 1349
$ws.Range("I61").Value = 1353.8889
$ws.Range("K61").Value = 1353.8889
$ws.Range("M61").Value = -1151.8889

$ws.Range("H94").Value = 55333.332
$ws.Range("J94").Value = 55333.332
$ws.Range("L94").Value = 55333.332
$ws.Range("N94").Value = -56685.332

$ws.Range("H113").Value = 1349
$ws.Range("I113").Value = 1353.8889
$ws.Range("K113").Value = 1353.8889
$ws.Range("M113").Value = 816.1111000000001

$ws.Range("H136").Value = 2847.4
$ws.Range("I136").Value = 2585.3572
$ws.Range("J136").Value = 3180.9092
$ws.Range("K136").Value = 7756.071599999999
$ws.Range("L136").Value = 9542.7276
$ws.Range("M136").Value = -5206.071599999999
$ws.Range("N136").Value = -14642.7276

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 25928.6
$ws.Range("J95").Value = 25928.6
$ws.Range("L95").Value = 25928.6
$ws.Range("N95").Value = -31420.6

$ws.Range("H107").Value = 391.72726
$ws.Range("I107").Value = 190.9
$ws.Range("K107").Value = 572.7
$ws.Range("M107").Value = 1347.3

$ws.Range("H135").Value = 44200
$ws.Range("J135").Value = 44200
$ws.Range("L135").Value = 44200
$ws.Range("N135").Value = -54340

$ws.Range("H136").Value = 1943.0714
$ws.Range("I136").Value = 1938.6923
$ws.Range("K136").Value = 5816.0769
$ws.Range("M136").Value = -3266.0769
